# #3518 - Update Personnel and Investment Item Template
#
# The "Section Org" code column (D) on the HR-Expense sheet used a mix of
# 6-digit codes (some stored as text, some as numbers). They are updated to
# the new 7-digit numeric Section Org codes used by the latest org mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HR-Expense")

# D10 kept an inconsistent "general" alignment style inherited from the old
# template; normalise it to match its siblings (left aligned) before writing
# the new value.
$ws.Range("D10").HorizontalAlignment = -4131

$ws.Range("D9").Value  = 1099999
$ws.Range("D10").Value = 1099999
$ws.Range("D11").Value = 2099999
$ws.Range("D12").Value = 2099999
$ws.Range("D13").Value = 3099999
$ws.Range("D14").Value = 3099999
$ws.Range("D15").Value = 4099999
$ws.Range("D16").Value = 4099999
$ws.Range("D17").Value = 5099999
$ws.Range("D18").Value = 5099999

# Drop the sheet-protection password (protection itself stays enabled).
$ws.Unprotect("cf2a")
$ws.Protect()

# Leave the cursor parked on F5, as in the saved file.
$ws.Activate()
$ws.Range("F5").Select()
